# edit.ps1 - Applies the "write up of statistical analysis" commit.
#
# 1. Removes the standalone "Results" heading paragraph that used to sit
#    right before "Plots" (the two headings get merged into one "Plots"
#    section at the top).
# 2. Appends a new "Statistical Analyses" section and a new "Results"
#    section (with write-up paragraphs) at the very end of the document,
#    just before the closing section properties.

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $exactText) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $exactText) {
            return $p
        }
    }
    return $null
}

# --- Step 1: remove the old standalone "Results" heading paragraph -------
$resultsHeading = Find-ParagraphByText $d "Results"
if ($resultsHeading -ne $null) {
    $resultsHeading.Range.Delete()
}

Write-Output ("Paragraph count after removing old Results heading: " + $d.Paragraphs.Count)

# --- Helpers for building the new content at the end of the document -----

# Appends a brand new empty paragraph at the very end of the document body
# (right before the sectPr) and returns it.
function New-EndParagraph($doc) {
    $lastIndex = $doc.Paragraphs.Count
    $lastPara = $doc.Paragraphs.Item($lastIndex)
    $lastPara.Range.InsertParagraphAfter()
    $newIndex = $doc.Paragraphs.Count
    return $doc.Paragraphs.Item($newIndex)
}

# Appends plain text (no special formatting) at the given (collapsed)
# position and returns the position right after the inserted text.
function Add-PlainRun($doc, $pos, $text) {
    $r = $doc.Range($pos, $pos)
    $r.InsertAfter($text)
    return $r.End
}

# Appends text at the given position and applies formatting to just the
# inserted span. $kind is one of "italic", "superscript", "subscript".
# Returns the new position right after the inserted text.
function Add-FormattedRun($doc, $pos, $text, $kind) {
    $start = $pos
    $r = $doc.Range($pos, $pos)
    $r.InsertAfter($text)
    $end = $r.End
    $span = $doc.Range($start, $end)
    if ($kind -eq "italic") {
        $span.Font.Italic = $true
    } elseif ($kind -eq "superscript") {
        $span.Font.Superscript = $true
    } elseif ($kind -eq "subscript") {
        $span.Font.Subscript = $true
    }
    return $end
}

Write-Output "helpers defined"

# Creates a heading paragraph (style "berschrift1") with the given text,
# wrapped in a bookmark of the given name.
function Add-Heading($doc, $text, $bookmarkName) {
    $p = New-EndParagraph $doc
    $p.Style = "berschrift1"
    $start = $p.Range.Start
    $pos = Add-PlainRun $doc $start $text
    $bmRange = $doc.Range($start, $pos)
    $doc.Bookmarks.Add($bookmarkName, $bmRange) | Out-Null
    return $p
}

# --- Step 2: append the new "Statistical Analyses" section ---------------
Add-Heading $d "Statistical Analyses" "statistical-analyses" | Out-Null

$p2 = New-EndParagraph $d
$p2.Style = "FirstParagraph"
$pos = $p2.Range.Start
$pos = Add-PlainRun $d $pos "The time course data from the eye-tracking task were analized using growth"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "curve analysis (GCA, Mirman, 2016)."
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "We downsampled the data to bins of 50 ms which were centered at the offset of"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "the first syllable of target items."
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "The time course of fixation ranged from 200 ms before target syllable offset"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "to 600 ms after."
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "The empirical logit transformation (Barr, 2008) was applied to the binary"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "responses (fixations to the target or the distractor)."
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "We modeled the time course using linear, quadratic, and cubic orthogonal"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "polynomials with fixed effects of group, lexical stress, and syllable"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "structure on all time terms."
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "For the group predictor M was set as the baseline, thus the IN and NIN"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "parameters described how the growth curve of the learners differed from that"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "of the native controls."
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "Lexical stress and syllable structure were sum coded such that parameter"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "estimates represent effect sizes of change from CV to CVC syllables and"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "paroxytone to oxytone stress."
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "All models included by-subject random effects on all time terms and the"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "syllable structure and lexical stress predictors, as well as by-item random"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "effects on all time terms."
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "Main effects and higher order interactions were assessed using nested model"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "comparisons."
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "The analysis was conducted in R (R Core Team, 2019) and the GCA models were"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "fit using lme4 (Bates, Mächler, Bolker, & Walker, 2009)."
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "Pairwise comparisons between learners groups were conducted using the R"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "package multcomp (Hothorn, Bretz, & Westfall, 2008)."

Write-Output ("Statistical Analyses paragraph length: " + $p2.Range.Text.Length)

# --- Step 3: append the new "Results" section -----------------------------
Add-Heading $d "Results" "results" | Out-Null

$p3 = New-EndParagraph $d
$p3.Style = "FirstParagraph"
$pos = $p3.Range.Start
$pos = Add-PlainRun $d $pos "The was a main effect of something on something"
$pos = Add-PlainRun $d $pos " "
$pos = Add-PlainRun $d $pos "("
$pos = Add-FormattedRun $d $pos "χ" "italic"
$pos = Add-FormattedRun $d $pos "2" "superscript"
$pos = Add-PlainRun $d $pos "(2) = 11,"
$pos = Add-PlainRun $d $pos " "
$pos = Add-FormattedRun $d $pos "p" "italic"
$pos = Add-PlainRun $d $pos " = .004)."

Write-Output ("Results paragraph 1 text: [" + $p3.Range.Text + "]")

$p4 = New-EndParagraph $d
$p4.Style = "Textkrper"
$pos = $p4.Range.Start
$pos = Add-PlainRun $d $pos "The effect was awesome ("
$pos = Add-PlainRun $d $pos "γ"
$pos = Add-FormattedRun $d $pos "00" "subscript"
$pos = Add-PlainRun $d $pos " = 1.18; SE = 0.22;"
$pos = Add-PlainRun $d $pos " "
$pos = Add-FormattedRun $d $pos "t" "italic"
$pos = Add-PlainRun $d $pos " = 5.36;"
$pos = Add-PlainRun $d $pos " "
$pos = Add-FormattedRun $d $pos "p" "italic"
$pos = Add-PlainRun $d $pos " < .001)."

Write-Output ("Results paragraph 2 text: [" + $p4.Range.Text + "]")

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
Write-Output ("Final bookmark count: " + $d.Bookmarks.Count)
